$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 160, pushing existing rows 160-192 down to 161-193.
$ws.Rows.Item(160).Insert()

$ws.Cells.Item(160, 1).Value = 6
$ws.Cells.Item(160, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(160, 3).Value = "Metropolitana"
$ws.Cells.Item(160, 4).Value = 44637
$ws.Cells.Item(160, 5).Value = 13
$ws.Cells.Item(160, 6).Value = 100112001
$ws.Cells.Item(160, 7).Value = "Berenjena"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 410
$ws.Cells.Item(160, 11).Value = 8000
$ws.Cells.Item(160, 12).Value = 9000
$ws.Cells.Item(160, 13).Value = 8561
$ws.Cells.Item(160, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(160, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(160, 16).Value = 171
$ws.Cells.Item(160, 17).Value = 50
$ws.Cells.Item(160, 18).Value = "Hortaliza"
